# Update "Pagos" (F), and "Inscrições homologadas" (H) counts, and
# "Inscritos" (E) counts on the "Inscricoes" worksheet per the latest
# enrollment/payment figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Cells.Item(2, 6).Value = 52
$ws.Cells.Item(2, 8).Value = 52
$ws.Cells.Item(3, 5).Value = 26
$ws.Cells.Item(3, 6).Value = 22
$ws.Cells.Item(3, 8).Value = 22
$ws.Cells.Item(4, 5).Value = 29
$ws.Cells.Item(4, 6).Value = 19
$ws.Cells.Item(4, 8).Value = 19
$ws.Cells.Item(5, 6).Value = 56
$ws.Cells.Item(5, 8).Value = 56
$ws.Cells.Item(6, 6).Value = 17
$ws.Cells.Item(6, 8).Value = 17
$ws.Cells.Item(7, 6).Value = 15
$ws.Cells.Item(7, 8).Value = 15
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(8, 8).Value = 3
$ws.Cells.Item(10, 5).Value = 246
$ws.Cells.Item(10, 6).Value = 132
$ws.Cells.Item(10, 8).Value = 132
$ws.Cells.Item(11, 5).Value = 182
$ws.Cells.Item(11, 6).Value = 116
$ws.Cells.Item(11, 8).Value = 116
$ws.Cells.Item(12, 5).Value = 271
$ws.Cells.Item(12, 6).Value = 167
$ws.Cells.Item(12, 8).Value = 167
$ws.Cells.Item(13, 5).Value = 86
$ws.Cells.Item(13, 6).Value = 44
$ws.Cells.Item(13, 8).Value = 44
$ws.Cells.Item(14, 5).Value = 74
$ws.Cells.Item(14, 6).Value = 39
$ws.Cells.Item(14, 8).Value = 39
$ws.Cells.Item(15, 6).Value = 39
$ws.Cells.Item(15, 8).Value = 39
$ws.Cells.Item(16, 5).Value = 110
$ws.Cells.Item(16, 6).Value = 70
$ws.Cells.Item(16, 8).Value = 70
$ws.Cells.Item(17, 5).Value = 53
$ws.Cells.Item(17, 6).Value = 29
$ws.Cells.Item(17, 8).Value = 29
$ws.Cells.Item(18, 6).Value = 20
$ws.Cells.Item(18, 8).Value = 20
$ws.Cells.Item(19, 6).Value = 5
$ws.Cells.Item(19, 8).Value = 5
$ws.Cells.Item(20, 6).Value = 25
$ws.Cells.Item(20, 8).Value = 25
$ws.Cells.Item(21, 5).Value = 83
$ws.Cells.Item(21, 6).Value = 51
$ws.Cells.Item(21, 8).Value = 51
$ws.Cells.Item(22, 6).Value = 62
$ws.Cells.Item(22, 8).Value = 62
$ws.Cells.Item(23, 6).Value = 63
$ws.Cells.Item(23, 8).Value = 63
$ws.Cells.Item(24, 5).Value = 122
$ws.Cells.Item(24, 6).Value = 67
$ws.Cells.Item(24, 8).Value = 67
$ws.Cells.Item(25, 5).Value = 126
$ws.Cells.Item(25, 6).Value = 67
$ws.Cells.Item(25, 8).Value = 67
$ws.Cells.Item(26, 6).Value = 51
$ws.Cells.Item(26, 8).Value = 51
$ws.Cells.Item(27, 5).Value = 176
$ws.Cells.Item(27, 6).Value = 102
$ws.Cells.Item(27, 8).Value = 102
$ws.Cells.Item(28, 6).Value = 37
$ws.Cells.Item(28, 8).Value = 37
$ws.Cells.Item(29, 6).Value = 72
$ws.Cells.Item(29, 8).Value = 72
$ws.Cells.Item(30, 5).Value = 127
$ws.Cells.Item(30, 6).Value = 77
$ws.Cells.Item(30, 8).Value = 77
$ws.Cells.Item(31, 5).Value = 46
$ws.Cells.Item(31, 6).Value = 25
$ws.Cells.Item(31, 8).Value = 25
$ws.Cells.Item(32, 5).Value = 118
$ws.Cells.Item(32, 6).Value = 67
$ws.Cells.Item(32, 8).Value = 67
$ws.Cells.Item(33, 6).Value = 93
$ws.Cells.Item(33, 8).Value = 93
$ws.Cells.Item(34, 6).Value = 85
$ws.Cells.Item(34, 8).Value = 85
$ws.Cells.Item(35, 5).Value = 90
$ws.Cells.Item(35, 6).Value = 58
$ws.Cells.Item(35, 8).Value = 58
$ws.Cells.Item(36, 5).Value = 33
$ws.Cells.Item(36, 6).Value = 23
$ws.Cells.Item(36, 8).Value = 23
$ws.Cells.Item(37, 5).Value = 89
$ws.Cells.Item(37, 6).Value = 53
$ws.Cells.Item(37, 8).Value = 53
$ws.Cells.Item(38, 5).Value = 56
$ws.Cells.Item(38, 6).Value = 42
$ws.Cells.Item(38, 8).Value = 42
$ws.Cells.Item(39, 5).Value = 121
$ws.Cells.Item(39, 6).Value = 63
$ws.Cells.Item(39, 8).Value = 63
$ws.Cells.Item(40, 5).Value = 159
$ws.Cells.Item(40, 6).Value = 81
$ws.Cells.Item(40, 8).Value = 81
$ws.Cells.Item(41, 5).Value = 220
$ws.Cells.Item(41, 6).Value = 107
$ws.Cells.Item(41, 8).Value = 107
$ws.Cells.Item(42, 5).Value = 205
$ws.Cells.Item(42, 6).Value = 117
$ws.Cells.Item(42, 8).Value = 117
$ws.Cells.Item(43, 5).Value = 64
$ws.Cells.Item(43, 6).Value = 34
$ws.Cells.Item(43, 8).Value = 34
$ws.Cells.Item(44, 5).Value = 172
$ws.Cells.Item(44, 6).Value = 100
$ws.Cells.Item(44, 8).Value = 100
$ws.Cells.Item(45, 6).Value = 42
$ws.Cells.Item(45, 8).Value = 42
$ws.Cells.Item(46, 6).Value = 96
$ws.Cells.Item(46, 8).Value = 96
$ws.Cells.Item(47, 5).Value = 257
$ws.Cells.Item(47, 6).Value = 140
$ws.Cells.Item(47, 8).Value = 140
$ws.Cells.Item(48, 5).Value = 123
$ws.Cells.Item(48, 6).Value = 55
$ws.Cells.Item(48, 8).Value = 55
$ws.Cells.Item(49, 5).Value = 141
$ws.Cells.Item(49, 6).Value = 72
$ws.Cells.Item(49, 8).Value = 72
$ws.Cells.Item(50, 5).Value = 114
$ws.Cells.Item(50, 6).Value = 55
$ws.Cells.Item(50, 8).Value = 55
$ws.Cells.Item(51, 5).Value = 120
$ws.Cells.Item(51, 6).Value = 59
$ws.Cells.Item(51, 8).Value = 59
$ws.Cells.Item(52, 5).Value = 12
